$d = $word.ActiveDocument

# Locate the "Visualisation : tableau" paragraph - the new predecessor
# for "Langages : r, python, matlab, c, c++" inside the COMPETENCES
# TECHNIQUES block (that paragraph currently sits further down, right
# after "MLOps : ...").
$visPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($visPara -eq $null -and $t -eq "Visualisation : tableau`r") {
        $visPara = $p
    }
}

# Re-create the "Langages" paragraph immediately before "Visualisation"
# (same paragraph formatting is inherited from the Visualisation
# paragraph, which carries the identical <w:spacing .../> pPr).
$visPara.Range.InsertBefore("Langages : r, python, matlab, c, c++`r")

# Remove the original "Langages" paragraph further down the list (the
# second occurrence of that exact text, since the first is the one we
# just inserted).
$seen = 0
$oldLangPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -eq "Langages : r, python, matlab, c, c++`r") {
        $seen = $seen + 1
        if ($seen -eq 2) {
            $oldLangPara = $p
        }
    }
}
$oldLangPara.Range.Delete()
